$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update date formatting from DD/MM/YYYY to DD-MM-YYYY for rows 3-21 (col A).
# Some day numbers are <= 12, so Excel's automatic type inference would
# otherwise mis-read "DD-MM-YYYY" as a US-style "MM-DD-YYYY" date and coerce
# it to a date serial number. Prefix those with a literal quote so Excel
# keeps them as plain text, exactly like typing '01-08-2022 into a cell.
$ws.Range("A3").Value = "28-07-2022"
$ws.Range("A4").Value = "'01-08-2022"
$ws.Range("A5").Value = "'04-08-2022"
$ws.Range("A6").Value = "'08-08-2022"
$ws.Range("A7").Value = "'11-08-2022"
$ws.Range("A8").Value = "15-08-2022"
$ws.Range("A9").Value = "18-08-2022"
$ws.Range("A10").Value = "22-08-2022"
$ws.Range("A11").Value = "25-08-2022"
$ws.Range("A12").Value = "29-08-2022"
$ws.Range("A13").Value = "'01-09-2022"
$ws.Range("A14").Value = "'05-09-2022"
$ws.Range("A15").Value = "'08-09-2022"
$ws.Range("A16").Value = "'12-09-2022"
$ws.Range("A17").Value = "15-09-2022"
$ws.Range("A18").Value = "19-09-2022"
$ws.Range("A19").Value = "22-09-2022"
$ws.Range("A20").Value = "26-09-2022"
$ws.Range("A21").Value = "29-09-2022"

# Update attendance counts for row 3 (28-07-2022): D3 0->1, G3 0->1
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# Update attendance counts for row 5 (04-08-2022): D5 0->1, E5 0->1, H5 1->0
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0
